$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B25").Value = 6474
$ws.Range("C25").Value = 1009
$ws.Range("D25").Value = 6015758
$ws.Range("E25").Value = 929.2181031819586
$ws.Range("F25").Value = 9.896452215243601
$ws.Range("G25").Value = 7.569296375266532
$ws.Range("H25").Value = 25.98852402779064
